# Add new "addons" column (I) to the Cigna info sheet with its value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "addons"
$ws.Range("I2").Value = "Dental/Maternity (Consultations, Scans and Delivery)/Optical Benefits/Wellness & Health Screening"

# Match the cursor position left behind in the saved workbook.
$ws.Range("H5").Select()
